$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = -102.235
$ws.Range("C22").Value = -102.235
$ws.Range("C22").Select()
